# Update KiCost/KiBot BoM-Costs workbook:
#  - bump the KiCad version string
#  - bump the "Created" timestamp and the KiCost/KiBot version string
#  - shrink a couple of column widths (D on BoM/DNF, E on Costs/Costs (DNF))

$wb = $excel.ActiveWorkbook

$newKicadVersion = "6.0.11+dfsg-1"
$newCreated      = "2024-02-02 17:09:04"
$newKicostKibot  = "KiCost" + [char]0x00AE + " v1.1.18 + KiBot v1.6.4"

# --- BoM sheet ---
$ws = $wb.Worksheets.Item("BoM")
$ws.Range("D6").Value = $newKicadVersion
$ws.Columns.Item(4).ColumnWidth = 20

# --- DNF sheet ---
$ws = $wb.Worksheets.Item("DNF")
$ws.Range("D6").Value = $newKicadVersion
$ws.Columns.Item(4).ColumnWidth = 20

# --- Costs sheet ---
$ws = $wb.Worksheets.Item("Costs")
$ws.Range("E6").Value = $newKicadVersion
$ws.Range("B24").Value = $newCreated
$ws.Range("A25").Value = $newKicostKibot
$ws.Columns.Item(5).ColumnWidth = 17

# --- Costs (DNF) sheet ---
$ws = $wb.Worksheets.Item("Costs (DNF)")
$ws.Range("E6").Value = $newKicadVersion
$ws.Range("B21").Value = $newCreated
$ws.Range("A22").Value = $newKicostKibot
$ws.Columns.Item(5).ColumnWidth = 17
